$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.05880547838015793
$ws.Range("C2").Value = 0.6749290290549438
$ws.Range("D2").Value = 1.036571427464883
$ws.Range("E2").Value = 1.018121519006883
$ws.Range("F2").Value = 1.040340913673013
$ws.Range("G2").Value = 22

$ws.Range("B3").Value = 0.1087012210528854
$ws.Range("C3").Value = 0.9180492773297835
$ws.Range("D3").Value = 1.469677934780255
$ws.Range("E3").Value = 1.212302740564524
$ws.Range("F3").Value = 1.2358333133304

$ws.Range("B4").Value = 0.09424136353897376
$ws.Range("C4").Value = 1.505548170551489
$ws.Range("D4").Value = 10.11276014086356
$ws.Range("E4").Value = 3.180056625417786
$ws.Range("F4").Value = 3.257157141062581
